$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-13 (columns A-T)
# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E=Ligand-expressing cells, F=Ligand detection rate, G=Ligand average expression value,
# H=Ligand total expression value, I=Ligand derived specificity of average expression value,
# J=Ligand derived specificity of total expression value, K=Receptor-expressing cells,
# L=Receptor detection rate, M=Receptor average expression value, N=Receptor total expression value,
# O=Receptor derived specificity of average expression value, P=Receptor derived specificity of total expression value,
# Q=Edge average expression weight, R=Edge total expression weight,
# S=Edge average expression derived specificity, T=Edge total expression derived specificity

$data = @(
  @("ECs","Ccl25","Ccr10","ECs",2,1,4.6982685,9.396537,0.1448422764790916,0.108319921689495,1,0.5,0.2651005,0.530201,0.2682690254597046,0.2682690254597046,1.24551332848425,4.982053313937,0.03885669635641101,0.02905887982951234),
  @("ECs","Ccl25","Ccr10","MuSCs",2,1,4.6982685,9.396537,0.1448422764790916,0.108319921689495,1,0.5,0.7230885,1.446177,0.7317309745402955,0.7317309745402955,3.39726392226225,13.589055689049,0.1059855801226806,0.07926104185998266),
  @("FAPs","Ccl25","Ccr10","ECs",3,1,5.277942,15.833826,0.1627129514638446,0.1825266896054461,1,0.5,0.2651005,0.530201,0.2682690254597046,0.2682690254597046,1.399185063171,8.395110379026002,0.04365084491887781,0.04896625714083903),
  @("FAPs","Ccl25","Ccr10","MuSCs",3,1,5.277942,15.833826,0.1627129514638446,0.1825266896054461,1,0.5,0.7230885,1.446177,0.7317309745402955,0.7317309745402955,3.816419163867,22.898514983202,0.1190621065449668,0.1335604324646071),
  @("Inflammatory-Mac","Ccl25","Ccr10","ECs",3,1,5.274651666666667,15.823955,0.162611514227898,0.1824129002437912,1,0.5,0.2651005,0.530201,0.2682690254597046,0.2682690254597046,1.398312794159167,8.389876764955002,0.04362363245044507,0.04893573097968017),
  @("Inflammatory-Mac","Ccl25","Ccr10","MuSCs",3,1,5.274651666666667,15.823955,0.162611514227898,0.1824129002437912,1,0.5,0.7230885,1.446177,0.7317309745402955,0.7317309745402955,3.8140399616725,22.884239770035,0.1189878817774529,0.133477169264111),
  @("MuSCs","Ccl25","Ccr10","ECs",2,1,5.865133,11.730266,0.1808153824270886,0.1352223158932856,1,0.5,0.2651005,0.530201,0.2682690254597046,0.2682690254597046,1.5548496908665,6.219398763466001,0.04850716643183887,0.03627595890509604),
  @("MuSCs","Ccl25","Ccr10","MuSCs",2,1,5.865133,11.730266,0.1808153824270886,0.1352223158932856,1,0.5,0.7230885,1.446177,0.7317309745402955,0.7317309745402955,4.2410102232705,16.964040893082,0.1323082159952498,0.09894635698818953),
  @("Neutrophils","Ccl25","Ccr10","ECs",3,1,5.310348,15.931044,0.1637119916020532,0.183647383979002,1,0.5,0.2651005,0.530201,0.2682690254597046,0.2682690254597046,1.407775909974,8.446655459844001,0.04391885644315017,0.04926690472827103),
  @("Neutrophils","Ccl25","Ccr10","MuSCs",3,1,5.310348,15.931044,0.1637119916020532,0.183647383979002,1,0.5,0.7230885,1.446177,0.7317309745402955,0.7317309745402955,3.839851569798,23.039109418788,0.1197931351589031,0.1343804792507309),
  @("Resolving-Mac","Ccl25","Ccr10","ECs",3,1,6.010791999999999,18.032376,0.1853058838000238,0.2078707885889801,1,0.5,0.2651005,0.530201,0.2682690254597046,0.2682690254597046,1.593463964596,9.560783787576,0.04971182885898165,0.05576529387630597),
  @("Resolving-Mac","Ccl25","Ccr10","MuSCs",3,1,6.010791999999999,18.032376,0.1853058838000238,0.2078707885889801,1,0.5,0.7230885,1.446177,0.7317309745402955,0.7317309745402955,4.346334571092,26.078007426552,0.1355940549410422,0.1521054947126741)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
}
